$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.969508
$ws.Range("H2").Value = 116.908524
$ws.Range("I2").Value = 0.688733638790647
$ws.Range("J2").Value = 0.688733638790647
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 405.24646
$ws.Range("N2").Value = 1215.73938
$ws.Range("O2").Value = 0.971171031955694
$ws.Range("P2").Value = 0.9711710319556939
$ws.Range("Q2").Value = 15792.25516494168
$ws.Range("R2").Value = 142130.2964844751
$ws.Range("S2").Value = 0.6688781587269128
$ws.Range("T2").Value = 0.6688781587269127
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.969508
$ws.Range("H3").Value = 116.908524
$ws.Range("I3").Value = 0.688733638790647
$ws.Range("J3").Value = 0.688733638790647
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.631177666666667
$ws.Range("N3").Value = 16.893533
$ws.Range("O3").Value = 0.01349508796612936
$ws.Range("P3").Value = 0.01349508796612936
$ws.Range("Q3").Value = 219.444223130588
$ws.Range("R3").Value = 1974.998008175292
$ws.Range("S3").Value = 0.009294521040712147
$ws.Range("T3").Value = 0.009294521040712147
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 38.969508
$ws.Range("H4").Value = 116.908524
$ws.Range("I4").Value = 0.688733638790647
$ws.Range("J4").Value = 0.688733638790647
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.958728333333333
$ws.Range("N4").Value = 5.876185
$ws.Range("O4").Value = 0.004694082255041018
$ws.Range("P4").Value = 0.004694082255041017
$ws.Range("Q4").Value = 76.33067945565999
$ws.Range("R4").Value = 686.97611510094
$ws.Range("S4").Value = 0.003232972352297006
$ws.Range("T4").Value = 0.003232972352297006
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 38.969508
$ws.Range("H5").Value = 116.908524
$ws.Range("I5").Value = 0.688733638790647
$ws.Range("J5").Value = 0.688733638790647
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.263573
$ws.Range("N5").Value = 9.790718999999999
$ws.Range("O5").Value = 0.007821135706583939
$ws.Range("P5").Value = 0.007821135706583937
$ws.Range("Q5").Value = 127.179834132084
$ws.Range("R5").Value = 1144.618507188756
$ws.Range("S5").Value = 0.005386679254671014
$ws.Range("T5").Value = 0.005386679254671013
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 38.969508
$ws.Range("H6").Value = 116.908524
$ws.Range("I6").Value = 0.688733638790647
$ws.Range("J6").Value = 0.688733638790647
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.176160333333333
$ws.Range("N6").Value = 3.528481
$ws.Range("O6").Value = 0.002818662116551706
$ws.Range("P6").Value = 0.002818662116551706
$ws.Range("Q6").Value = 45.83438951911599
$ws.Range("R6").Value = 412.509505672044
$ws.Range("S6").Value = 0.001941307416054003
$ws.Range("T6").Value = 0.001941307416054003
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.685730333333334
$ws.Range("H7").Value = 29.057191
$ws.Range("I7").Value = 0.1711822560557247
$ws.Range("J7").Value = 0.1711822560557247
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 405.24646
$ws.Range("N7").Value = 1215.73938
$ws.Range("O7").Value = 0.971171031955694
$ws.Range("P7").Value = 0.9711710319556939
$ws.Range("Q7").Value = 3925.107930097954
$ws.Range("R7").Value = 35325.97137088158
$ws.Range("S7").Value = 0.166247248266142
$ws.Range("T7").Value = 0.166247248266142
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.685730333333334
$ws.Range("H8").Value = 29.057191
$ws.Range("I8").Value = 0.1711822560557247
$ws.Range("J8").Value = 0.1711822560557247
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.631177666666667
$ws.Range("N8").Value = 16.893533
$ws.Range("O8").Value = 0.01349508796612936
$ws.Range("P8").Value = 0.01349508796612936
$ws.Range("Q8").Value = 54.54206833842256
$ws.Range("R8").Value = 490.8786150458031
$ws.Range("S8").Value = 0.002310119603712486
$ws.Range("T8").Value = 0.002310119603712486
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.685730333333334
$ws.Range("H9").Value = 29.057191
$ws.Range("I9").Value = 0.1711822560557247
$ws.Range("J9").Value = 0.1711822560557247
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.958728333333333
$ws.Range("N9").Value = 5.876185
$ws.Range("O9").Value = 0.004694082255041018
$ws.Range("P9").Value = 0.004694082255041017
$ws.Range("Q9").Value = 18.97171443292611
$ws.Range("R9").Value = 170.745429896335
$ws.Range("S9").Value = 0.0008035435905290654
$ws.Range("T9").Value = 0.0008035435905290652
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.685730333333334
$ws.Range("H10").Value = 29.057191
$ws.Range("I10").Value = 0.1711822560557247
$ws.Range("J10").Value = 0.1711822560557247
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.263573
$ws.Range("N10").Value = 9.790718999999999
$ws.Range("O10").Value = 0.007821135706583939
$ws.Range("P10").Value = 0.007821135706583937
$ws.Range("Q10").Value = 31.61008800114766
$ws.Range("R10").Value = 284.490792010329
$ws.Range("S10").Value = 0.001338839655171023
$ws.Range("T10").Value = 0.001338839655171023
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 9.685730333333334
$ws.Range("H11").Value = 29.057191
$ws.Range("I11").Value = 0.1711822560557247
$ws.Range("J11").Value = 0.1711822560557247
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.176160333333333
$ws.Range("N11").Value = 3.528481
$ws.Range("O11").Value = 0.002818662116551706
$ws.Range("P11").Value = 0.002818662116551706
$ws.Range("Q11").Value = 11.39197181743011
$ws.Range("R11").Value = 102.527746356871
$ws.Range("S11").Value = 0.0004825049401701252
$ws.Range("T11").Value = 0.0004825049401701252
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6067633333333333
$ws.Range("H12").Value = 1.82029
$ws.Range("I12").Value = 0.01072372580252768
$ws.Range("J12").Value = 0.01072372580252768
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 405.24646
$ws.Range("N12").Value = 1215.73938
$ws.Range("O12").Value = 0.971171031955694
$ws.Range("P12").Value = 0.9711710319556939
$ws.Range("Q12").Value = 245.8886928911333
$ws.Range("R12").Value = 2212.9982360202
$ws.Range("S12").Value = 0.01041457185405071
$ws.Range("T12").Value = 0.01041457185405071
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6067633333333333
$ws.Range("H13").Value = 1.82029
$ws.Range("I13").Value = 0.01072372580252768
$ws.Range("J13").Value = 0.01072372580252768
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.631177666666667
$ws.Range("N13").Value = 16.893533
$ws.Range("O13").Value = 0.01349508796612936
$ws.Range("P13").Value = 0.01349508796612936
$ws.Range("Q13").Value = 3.416792131618889
$ws.Range("R13").Value = 30.75112918457
$ws.Range("S13").Value = 0.0001447176230297623
$ws.Range("T13").Value = 0.0001447176230297623
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.6067633333333333
$ws.Range("H14").Value = 1.82029
$ws.Range("I14").Value = 0.01072372580252768
$ws.Range("J14").Value = 0.01072372580252768
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.958728333333333
$ws.Range("N14").Value = 5.876185
$ws.Range("O14").Value = 0.004694082255041018
$ws.Range("P14").Value = 0.004694082255041017
$ws.Range("Q14").Value = 1.188484532627778
$ws.Range("R14").Value = 10.69636079365
$ws.Range("S14").Value = 0.0000503380509975707
$ws.Range("T14").Value = 0.00005033805099757069
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.6067633333333333
$ws.Range("H15").Value = 1.82029
$ws.Range("I15").Value = 0.01072372580252768
$ws.Range("J15").Value = 0.01072372580252768
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.263573
$ws.Range("N15").Value = 9.790718999999999
$ws.Range("O15").Value = 0.007821135706583939
$ws.Range("P15").Value = 0.007821135706583937
$ws.Range("Q15").Value = 1.980216432056666
$ws.Range("R15").Value = 17.82194788851
$ws.Range("S15").Value = 0.00008387171478176477
$ws.Range("T15").Value = 0.00008387171478176475
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.6067633333333333
$ws.Range("H16").Value = 1.82029
$ws.Range("I16").Value = 0.01072372580252768
$ws.Range("J16").Value = 0.01072372580252768
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.176160333333333
$ws.Range("N16").Value = 3.528481
$ws.Range("O16").Value = 0.002818662116551706
$ws.Range("P16").Value = 0.002818662116551706
$ws.Range("Q16").Value = 0.7136509643877778
$ws.Range("R16").Value = 6.422858679489999
$ws.Range("S16").Value = 0.00003022655966787283
$ws.Range("T16").Value = 0.00003022655966787282
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4787493333333333
$ws.Range("H17").Value = 1.436248
$ws.Range("I17").Value = 0.008461250535040448
$ws.Range("J17").Value = 0.008461250535040449
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 405.24646
$ws.Range("N17").Value = 1215.73938
$ws.Range("O17").Value = 0.971171031955694
$ws.Range("P17").Value = 0.9711710319556939
$ws.Range("Q17").Value = 194.0114725606933
$ws.Range("R17").Value = 1746.10325304624
$ws.Range("S17").Value = 0.0082173214137509
$ws.Range("T17").Value = 0.0082173214137509
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.4787493333333333
$ws.Range("H18").Value = 1.436248
$ws.Range("I18").Value = 0.008461250535040448
$ws.Range("J18").Value = 0.008461250535040449
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 5.631177666666667
$ws.Range("N18").Value = 16.893533
$ws.Range("O18").Value = 0.01349508796612936
$ws.Range("P18").Value = 0.01349508796612936
$ws.Range("Q18").Value = 2.695922553798222
$ws.Range("R18").Value = 24.263302984184
$ws.Range("S18").Value = 0.00011418532027383
$ws.Range("T18").Value = 0.00011418532027383
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.4787493333333333
$ws.Range("H19").Value = 1.436248
$ws.Range("I19").Value = 0.008461250535040448
$ws.Range("J19").Value = 0.008461250535040449
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 1.958728333333333
$ws.Range("N19").Value = 5.876185
$ws.Range("O19").Value = 0.004694082255041018
$ws.Range("P19").Value = 0.004694082255041017
$ws.Range("Q19").Value = 0.9377398837644443
$ws.Range("R19").Value = 8.439658953879999
$ws.Range("S19").Value = 0.00003971780599198969
$ws.Range("T19").Value = 0.00003971780599198969
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 0.4787493333333333
$ws.Range("H20").Value = 1.436248
$ws.Range("I20").Value = 0.008461250535040448
$ws.Range("J20").Value = 0.008461250535040449
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 3.263573
$ws.Range("N20").Value = 9.790718999999999
$ws.Range("O20").Value = 0.007821135706583939
$ws.Range("P20").Value = 0.007821135706583937
$ws.Range("Q20").Value = 1.562433398034666
$ws.Range("R20").Value = 14.061900582312
$ws.Range("S20").Value = 0.0000661765886819573
$ws.Range("T20").Value = 0.0000661765886819573
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 0.4787493333333333
$ws.Range("H21").Value = 1.436248
$ws.Range("I21").Value = 0.008461250535040448
$ws.Range("J21").Value = 0.008461250535040449
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 1.176160333333333
$ws.Range("N21").Value = 3.528481
$ws.Range("O21").Value = 0.002818662116551706
$ws.Range("P21").Value = 0.002818662116551706
$ws.Range("Q21").Value = 0.5630859754764443
$ws.Range("R21").Value = 5.067773779287999
$ws.Range("S21").Value = 0.00002384940634177137
$ws.Range("T21").Value = 0.00002384940634177137
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 6.840640999999999
$ws.Range("H22").Value = 20.521923
$ws.Range("I22").Value = 0.1208991288160602
$ws.Range("J22").Value = 0.1208991288160602
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 405.24646
$ws.Range("N22").Value = 1215.73938
$ws.Range("O22").Value = 0.971171031955694
$ws.Range("P22").Value = 0.9711710319556939
$ws.Range("Q22").Value = 2772.14554938086
$ws.Range("R22").Value = 24949.30994442774
$ws.Range("S22").Value = 0.1174137316948376
$ws.Range("T22").Value = 0.1174137316948376
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 6.840640999999999
$ws.Range("H23").Value = 20.521923
$ws.Range("I23").Value = 0.1208991288160602
$ws.Range("J23").Value = 0.1208991288160602
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 5.631177666666667
$ws.Range("N23").Value = 16.893533
$ws.Range("O23").Value = 0.01349508796612936
$ws.Range("P23").Value = 0.01349508796612936
$ws.Range("Q23").Value = 38.52086482488433
$ws.Range("R23").Value = 346.687783423959
$ws.Range("S23").Value = 0.001631544378401138
$ws.Range("T23").Value = 0.001631544378401138
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 6.840640999999999
$ws.Range("H24").Value = 20.521923
$ws.Range("I24").Value = 0.1208991288160602
$ws.Range("J24").Value = 0.1208991288160602
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 1.958728333333333
$ws.Range("N24").Value = 5.876185
$ws.Range("O24").Value = 0.004694082255041018
$ws.Range("P24").Value = 0.004694082255041017
$ws.Range("Q24").Value = 13.39895734486166
$ws.Range("R24").Value = 120.590616103755
$ws.Range("S24").Value = 0.0005675104552253866
$ws.Range("T24").Value = 0.0005675104552253865
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 6.840640999999999
$ws.Range("H25").Value = 20.521923
$ws.Range("I25").Value = 0.1208991288160602
$ws.Range("J25").Value = 0.1208991288160602
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 3.263573
$ws.Range("N25").Value = 9.790718999999999
$ws.Range("O25").Value = 0.007821135706583939
$ws.Range("P25").Value = 0.007821135706583937
$ws.Range("Q25").Value = 22.32493127029299
$ws.Range("R25").Value = 200.9243814326369
$ws.Range("S25").Value = 0.0009455684932781798
$ws.Range("T25").Value = 0.0009455684932781796
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 6.840640999999999
$ws.Range("H26").Value = 20.521923
$ws.Range("I26").Value = 0.1208991288160602
$ws.Range("J26").Value = 0.1208991288160602
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 1.176160333333333
$ws.Range("N26").Value = 3.528481
$ws.Range("O26").Value = 0.002818662116551706
$ws.Range("P26").Value = 0.002818662116551706
$ws.Range("Q26").Value = 8.045690598773664
$ws.Range("R26").Value = 72.41121538896299
$ws.Range("S26").Value = 0.0003407737943179337
$ws.Range("T26").Value = 0.0003407737943179336
